$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) from 2023-09-12 (45181) to 2023-09-13 (45182)
# for all data rows currently holding the old value.
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value = 45182
    }
}
